# Daily attendance processing - 2026-01-11 12:51:48
# Reorders the comma-separated "Recorded By" (column G) entries on the
# session analysis sheet so the human recorder's email/name is listed
# before the automated "System" marker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact-match replacements applied to column G ("Recorded By") values.
$replacements = @{
    'System, dnasr281@gmail.com'            = 'dnasr281@gmail.com, System'
    'System, system, backup@backdoor.com'   = 'system, System, backup@backdoor.com'
    'admin@admin.com, dnasr281@gmail.com'   = 'dnasr281@gmail.com, admin@admin.com'
}

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

$updatedCount = 0
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2
    if ($current -and $replacements.ContainsKey($current)) {
        $cell.Value = $replacements[$current]
        $updatedCount++
    }
}

Write-Host "Updated $updatedCount Recorded By cells"
